$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '66.896.08'
$ws.Cells.Item(2, 5).Value = '  -0.36%  '
$ws.Cells.Item(3, 4).Value = '3.084.90'
$ws.Cells.Item(3, 5).Value = '  +0.05%  '
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '571.07'
$ws.Cells.Item(5, 5).Value = '  -1.06%  '
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '176.24'
$ws.Cells.Item(6, 5).Value = '  +4.15%  '
$ws.Cells.Item(7, 5).Value = '  +0.09%  '
$ws.Cells.Item(8, 4).Value = '3.083.49'
$ws.Cells.Item(8, 5).Value = '  +0.10%  '
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.512'
$ws.Cells.Item(9, 5).Value = '  -0.53%  '
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '6.43'
$ws.Cells.Item(10, 5).Value = '  +0.33%  '
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.151'
$ws.Cells.Item(11, 5).Value = '  +0.07%  '
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '0.466'
$ws.Cells.Item(12, 5).Value = '  -1.29%  '
$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '0.0000240'
$ws.Cells.Item(13, 5).Value = '  -0.73%  '
$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '35.84'
$ws.Cells.Item(14, 5).Value = '  -0.93%  '
$ws.Cells.Item(15, 5).Value = '  +0.74%  '
$ws.Cells.Item(16, 4).Value = '3.594.57'
$ws.Cells.Item(16, 5).Value = '  -0.04%  '
$ws.Cells.Item(17, 4).Value = '66.853.63'
$ws.Cells.Item(17, 5).Value = '  -0.24%  '
$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '6.99'
$ws.Cells.Item(18, 5).Value = '  -0.49%  '
$ws.Cells.Item(19, 4).Value = '3.088.55'
$ws.Cells.Item(19, 5).Value = '  +0.09%  '
$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '16.48'
$ws.Cells.Item(20, 5).Value = '  -0.43%  '
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '485.18'
$ws.Cells.Item(21, 5).Value = '  -1.22%  '
$ws.Cells.Item(22, 5).Value = '  +0.00%  '
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '0.683'
$ws.Cells.Item(23, 5).Value = '  -0.77%  '
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '83.33'
$ws.Cells.Item(24, 5).Value = '  +0.52%  '
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '12.71'
$ws.Cells.Item(25, 5).Value = '  -1.67%  '
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '2.23'
$ws.Cells.Item(26, 5).Value = '  +0.15%  '
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '10.26'
$ws.Cells.Item(27, 5).Value = '  +0.19%  '
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.999'
$ws.Cells.Item(28, 5).Value = '  -0.07%  '
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '7.83'
$ws.Cells.Item(29, 5).Value = '  -0.61%  '
$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '2.28'
$ws.Cells.Item(30, 5).Value = '  -0.74%  '
$ws.Cells.Item(31, 5).Value = '  -1.85%  '
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '27.92'
$ws.Cells.Item(32, 5).Value = '  +0.09%  '
$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '0.111'
$ws.Cells.Item(33, 5).Value = '  -0.33%  '
$ws.Cells.Item(34, 4).Value = '0.0₃0935'
$ws.Cells.Item(34, 5).Value = '  +2.69%  '
$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '1.00'
$ws.Cells.Item(35, 5).Value = '  +0.31%  '
$ws.Cells.Item(36, 2).Value = 'Mantle'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.943'
$ws.Cells.Item(36, 5).Value = '  -1.33%  '
$ws.Cells.Item(37, 2).Value = 'Filecoin'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '5.54'
$ws.Cells.Item(37, 5).Value = '  -2.97%  '
$ws.Cells.Item(38, 2).Value = 'Arweave'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '46.81'
$ws.Cells.Item(38, 5).Value = '  -0.36%  '
$ws.Cells.Item(39, 5).Value = '  +2.24%  '
$ws.Cells.Item(40, 2).Value = 'OKB'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '48.91'
$ws.Cells.Item(40, 5).Value = '  -0.62%  '
$ws.Cells.Item(41, 2).Value = 'Stacks'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.99'
$ws.Cells.Item(41, 5).Value = '  -0.17%  '
$ws.Cells.Item(42, 5).Value = '  -0.36%  '
$ws.Cells.Item(43, 2).Value = 'Cosmos'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '8.21'
$ws.Cells.Item(43, 5).Value = '  -1.37%  '
$ws.Cells.Item(44, 2).Value = 'dogwifhat'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '2.69'
$ws.Cells.Item(44, 5).Value = '  +9.30%  '
$ws.Cells.Item(45, 4).Value = '2.801.78'
$ws.Cells.Item(45, 5).Value = '  +0.93%  '
$ws.Cells.Item(46, 4).NumberFormat = '@'
$ws.Cells.Item(46, 4).Value = '368.94'
$ws.Cells.Item(46, 5).Value = '  -0.76%  '
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '0.0342'
$ws.Cells.Item(47, 5).Value = '  -0.80%  '
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '134.59'
$ws.Cells.Item(48, 5).Value = '  -0.93%  '
$ws.Cells.Item(49, 5).Value = '  +0.03%  '
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '25.68'
$ws.Cells.Item(50, 5).Value = '  +4.69%  '
$ws.Cells.Item(51, 5).Value = '  +6.47%  '
